# This script reproduces, via the Excel COM object model, the edits described
# by the target diff:
#   - K3:K10 on Sheet1 get consolidated into a single shared formula group
#     (same formula text/values as before, now using t="shared")
#   - A new worksheet "Sheet2" is added after "Sheet1", populated with a small
#     RSA-notes table, and becomes the active/selected sheet
#   - Selection on Sheet1 moves to V17, selection on Sheet2 is H7

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1: turn the K3:K10 formulas into one shared-formula group ---
# (Re-applying the identical formula across the whole range makes Excel
# store it as t="shared" instead of one literal formula per cell.)
$ws1.Range("K3:K10").Formula = '=MID($R$2,$I3,1)*1'

# --- Add Sheet2 right after Sheet1 ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = "p,q"
$ws2.Range("B1").Value = "pierwsze"

$ws2.Range("A2").Value = "256<="
$ws2.Range("B2").Value = "p,q"
$ws2.Range("C2").Value = "<=512"

$ws2.Range("A3").Value = "n="
$ws2.Range("B3").Value = "p*q"

$ws2.Range("A4").Value = "Fi"
$ws2.Range("B4").Value = "(p-1)(q-1)"

$ws2.Range("A5").Value = "e="
$ws2.Range("B5").Value = 65537

$ws2.Range("A6").Value = "NWD(e, Fi)=1 ???"

# Note: B7/A7 are populated before C6 so that the new shared-string entries
# end up in the same append order as the target workbook.
$ws2.Range("B7").Value = "e*d mod Fi=1"
$ws2.Range("A7").Value = "d taka że"

$ws2.Range("C6").Value = "to szukamy d:"

# --- Selections & active sheet ---
$ws1.Range("V17").Select() | Out-Null
$ws2.Range("H7").Select() | Out-Null
$ws2.Activate() | Out-Null
